$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so Excel does not
# auto-convert them to Number cells (which would lose exact formatting,
# e.g. trailing zeros, and introduce floating point drift).
$ws.Range("D2").Value = "27.243.30"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.851.92"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.91"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4602"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3708"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07295"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8862"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.08"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07812"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "1.810.21"
$ws.Range("E13").Value = "  -7.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.384"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.522"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.38"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008932"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "27.279.14"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.107"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "2.045.41"
$ws.Range("E24").Value = "  -8.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.929"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.90"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.69"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.059"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08808"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.094"
$ws.Range("E32").Value = "  +4.40%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7699"
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.499"
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.759"
$ws.Range("E36").Value = "  +12.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05262"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.948"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.067"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5121"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1633"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.387"
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4796"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.30"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.25"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.641"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06208"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  +0.68%  "
